$wb = $excel.ActiveWorkbook

$logSheet = $wb.Worksheets.Item("Log")
$todoSheet = $wb.Worksheets.Item("To Do")

# --- Log sheet updates ---
# Row 1: text content unchanged
$logSheet.Range("B1").Value = "Replaced cadastre feature with region_id and uploaded to the repo"

# New "To Do" item introduced before the others
$todoSheet.Range("A2").Value = "Implement prediction with streamlit"

# Row 2 of the Log sheet gains a date and a note (previously blank besides the date style)
$logSheet.Range("A2").Value = "2022-03-19"
$logSheet.Range("A2").NumberFormat = $logSheet.Range("A1").NumberFormat
$logSheet.Range("B2").Value = "Read features from Webapp to dataframe and concatenated two dataframes labelling them 1 for training and 0 for validation"

# Another new "To Do" item, inserted at the very top
$todoSheet.Range("A1").Value = "Process all data n concatenated dataframe to get a normalized result"

# Remaining "To Do" items shift down to make room for the two new ones above
$todoSheet.Range("A3").Value = "Work with outliers of the numerical features (pricePerOne, estimatedPrice, etc)"
$todoSheet.Range("A4").Value = "Think what to do with owner edrpous"

# Widen column B on the Log sheet
$logSheet.Columns.Item(2).ColumnWidth = 108.6

# Update selections to match the saved view state
$logSheet.Activate()
$logSheet.Range("B6").Select()
$todoSheet.Activate()
$todoSheet.Range("F13").Select()
